$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (A21 s=9,10,12,4 ht=13.5 thickBot=1) has the exact formatting
# pattern we want for the new row. Copy it, then insert a copy of it at
# row 19 so all formatting (styles + row height) carries over, then
# overwrite the values with the Change Healthcare data.
$ws.Rows.Item(21).Copy()
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = "Change Healthcare"
$ws.Range("B19").Value = 7100
$ws.Range("C19").Value = "Customer Support"
$ws.Range("D19").Value = "change_Anrios"
